# Updates crypto price (D) and 1h volume change (E) cells to match latest
# scraped values, keeping cells as plain text (matching the original
# inlineStr cell type) with no style/number-format changes applied.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "66.238.71"
$ws.Range("D2").Style = "Normal"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.193.50"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +1.45%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "597.32"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +4.25%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "153.67"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +3.13%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.192.06"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +1.59%  "
$ws.Range("E9").Value = "  +2.01%  "
$ws.Range("E10").Value = "  +0.66%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.11"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.20%  "
$ws.Range("E12").Value = "  +3.82%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000269"
$ws.Range("D13").Style = "Normal"
$ws.Range("E14").Value = "  +5.37%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.717.95"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.49%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "66.191.15"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.89%  "
$ws.Range("E17").Value = "  +5.00%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.196.53"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.32%  "
$ws.Range("E19").Value = "  +0.64%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "510.17"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.21%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "15.31"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +3.75%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.737"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +3.26%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "8.00"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +3.93%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "15.06"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.81%  "
$ws.Range("E25").Value = "  +1.03%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.999"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.02%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.26"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +5.45%  "
$ws.Range("E28").Value = "  +3.18%  "
$ws.Range("E29").Value = "  +6.43%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "6.98"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +13.51%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.89"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +3.79%  "
$ws.Range("E32").Value = "  +2.58%  "
$ws.Range("E33").Value = "  +3.21%  "
$ws.Range("E34").Value = "  +0.05%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "6.52"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.56%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "54.82"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.24%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.0899"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.53%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "484.32"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +4.65%  "
$ws.Range("E39").Value = "  -0.15%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "8.83"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +2.44%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.90"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -2.63%  "
$ws.Range("E42").Value = "  +4.92%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.298"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +6.29%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0₃0647"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +11.61%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.919.28"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -4.22%  "
$ws.Range("E46").Value = "  -0.39%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "28.42"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.44%  "
$ws.Range("E48").Value = "  +0.00%  "
$ws.Range("E49").Value = "  +2.39%  "
$ws.Range("E50").Value = "  +3.63%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.58"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +5.88%  "
